$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("unit" / UNIT-ID), shifting the
# existing poll_period..topic columns from D:N to E:O.
$ws.Columns("D").Insert()

# Match the new column's display width to its left neighbour (column C),
# mirroring Excel's "insert copies formatting from the column to the left".
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth

# Row 2 (field-name header row): new field is "unit"
$ws.Range("D2").Value = "unit"

# Row 1 (description header row): rich text "UNIT-ID (1-127) " + bold "required"
$ws.Range("D1").Value = "UNIT-ID (1-127) required"
$reqChars = $ws.Range("D1").Characters(17, 8)
$reqChars.Font.Bold = $true

# Data rows: UNIT-ID value of 1 for each device row
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1

# Cursor ends on the newly added column's header cell
$ws.Range("D1").Select()
